$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E51').NumberFormat = "@"

$ws.Range("D2").Value = '68.267.37'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '3.347.46'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '582.91'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '177.33'
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = '0.183'
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("D11").Value = '48.03'
$ws.Range("E11").Value = '  +5.32%  '
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '690.85'
$ws.Range("E13").Value = '  +3.68%  '
$ws.Range("D14").Value = '3.885.82'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").Value = '8.43'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '68.259.58'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '0.120'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").Value = '3.346.03'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '17.45'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '11.19'
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").Value = '5.45'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '16.97'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = '100.01'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("D28").Value = '32.99'
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").Value = '8.49'
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").Value = '6.93'
$ws.Range("E30").Value = '  -6.91%  '
$ws.Range("D31").Value = '557.11'
$ws.Range("E31").Value = '  -5.99%  '
$ws.Range("D32").Value = '11.06'
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").Value = '0.106'
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("D34").Value = '57.64'
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '3.699.04'
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = '3.28'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("E38").Value = '  +3.32%  '
$ws.Range("D39").Value = '34.65'
$ws.Range("E39").Value = '  +4.03%  '
$ws.Range("E40").Value = '  +2.04%  '
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").Value = '0.0₃0672'
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").Value = '3.27'
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '0.0410'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("D46").Value = '2.65'
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '1.34'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").Value = '130.54'
$ws.Range("E50").Value = '  +2.73%  '
$ws.Range("E51").Value = '  -0.35%  '
